# "se añade parrafo del lobo"
#
# 1. Extend the first paragraph ("... del camino") with an extra run
#    describing the meeting with the woodcutter.
# 2. Insert a brand-new paragraph right after it describing the wolf,
#    including a misspelled word ("acimpañar") wrapped in the usual
#    w:proofErr spellStart/spellEnd markers.
# 3. Insert two genuinely empty paragraphs after that new paragraph.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Append the new sentence to the end of paragraph 1 -------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertAfter(" en el bosque se encontró con un leñador el cual le dijo que si la acompañaba pero ella le dijo que no que estaba muy cerca a la casa de la abuela")

# --- 2. Insert the new "lobo" paragraph right after paragraph 1 -------
$p1.Range.InsertParagraphAfter()
$loboPara = $d.Paragraphs.Item(2)

$loboXml = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">De repente apareció el lobo feroz que quería comerse a caperucita </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">y se le acerco para decirle que si la podía </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>acimpañar</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '</w:p>'

$loboPara.Range.InsertXML($loboXml)

# --- 3. Two empty paragraphs right after the "lobo" paragraph ---------
$loboPara = $d.Paragraphs.Item(2)
$loboPara.Range.InsertParagraphAfter()
$empty1 = $d.Paragraphs.Item(3)
$empty1.Range.InsertXML('<w:p ' + $wNs + '/>')

$empty1 = $d.Paragraphs.Item(3)
$empty1.Range.InsertParagraphAfter()
$empty2 = $d.Paragraphs.Item(4)
$empty2.Range.InsertXML('<w:p ' + $wNs + '/>')
